$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New survey rows for groups 112-117 (script codes 04,05,06,07,08,01), appended
# after the existing data (which ends at row 100 / group 111).
# Columns: RowNum, Group(A), Script(B, text), Event(C, text),
#          Disruption_Factor(D), Confident_Factor(E), Dispersion_School(F)
$rows = @(
    @(101, 112, "04", "head on table", 3, 5, 7),
    @(102, 112, "04", "drumming", 7, 3, 3),
    @(103, 112, "04", "snipping", 5, 9, 7),
    @(104, 112, "04", "clicking pen", 9, 3, 8),
    @(105, 112, "04", "heckling", 9, 4, 9),
    @(106, 112, "04", "locking at phone", 4, 6, 5),
    @(107, 112, "04", "whispering", 7, 6, 10),
    @(108, 112, "04", "drawing", 3, 8, 9),
    @(109, 112, "04", "chatting", 9, 6, 7),
    @(110, 113, "05", "drumming", $null, $null, $null),
    @(111, 113, "05", "clicking pen", 4, 10, 8),
    @(112, 113, "05", "head on table", 2, 10, 8),
    @(113, 113, "05", "locking at phone", 5, 8, 7),
    @(114, 113, "05", "snipping", 1, 10, 9),
    @(115, 113, "05", "drawing", 4, 9, 7),
    @(116, 113, "05", "heckling", 6, 10, 7),
    @(117, 113, "05", "whispering", 4, 10, 10),
    @(118, 113, "05", "chatting", 7, 8, 7),
    @(119, 114, "06", "clicking pen", 8, 10, 8),
    @(120, 114, "06", "locking at phone", 2, 9, 10),
    @(121, 114, "06", "drumming", 8, 9, 5),
    @(122, 114, "06", "drawing", 0, 10, 6),
    @(123, 114, "06", "head on table", 2, 9, 5),
    @(124, 114, "06", "whispering", 7, 9, 10),
    @(125, 114, "06", "snipping", 8, 10, 6),
    @(126, 114, "06", "heckling", 9, 6, 7),
    @(127, 114, "06", "chatting", 9, 9, 9),
    @(128, 115, "07", "locking at phone", $null, $null, $null),
    @(129, 115, "07", "drawing", 3, 8, 6),
    @(130, 115, "07", "clicking pen", 8, 6, 4),
    @(131, 115, "07", "whispering", 7, 6, 7),
    @(132, 115, "07", "drumming", 5, 8, 4),
    @(133, 115, "07", "heckling", 10, 1, 2),
    @(134, 115, "07", "head on table", 6, 8, 5),
    @(135, 115, "07", "snipping", 9, 8, 8),
    @(136, 115, "07", "chatting", 10, 0, 1),
    @(137, 116, "08", "drawing", 0, 10, 8),
    @(138, 116, "08", "whispering", 8, 7, 7),
    @(139, 116, "08", "locking at phone", 5, 6, 8),
    @(140, 116, "08", "heckling", 10, 3, 3),
    @(141, 116, "08", "clicking pen", 8, 7, 5),
    @(142, 116, "08", "snipping", 0, 10, 5),
    @(143, 116, "08", "drumming", 9, 10, 5),
    @(144, 116, "08", "head on table", 8, 10, 7),
    @(145, 116, "08", "chatting", 10, 8, 3),
    @(146, 117, "01", "whispering", 3, 4, 10),
    @(147, 117, "01", "heckling", 7, 3, 3),
    @(148, 117, "01", "drawing", 1, 5, 9),
    @(149, 117, "01", "snipping", 4, 6, 5),
    @(150, 117, "01", "locking at phone", 4, 6, 9),
    @(151, 117, "01", "head on table", 1, 4, 7),
    @(152, 117, "01", "clicking pen", 5, 6, 6),
    @(153, 117, "01", "drumming", 5, 5, 7),
    @(154, 117, "01", "chatting", 7, 4, 7)
)

# Source rows whose formatting should be replicated:
#  - row 99 carries the plain interior styling (style 1 / 2 for the text column)
#  - row 100 carries the "last row of group" styling with the bottom border
#    (style 3 / 4 for the text column)
$normalFormatRow = 99
$lastFormatRow = 100

for ($i = 0; $i -lt $rows.Count; $i++) {
    $d = $rows[$i]
    $r = $d[0]
    $isLastOfGroup = ($i -eq ($rows.Count - 1)) -or ($rows[$i + 1][1] -ne $d[1])

    if ($isLastOfGroup) {
        $ws.Range("A$lastFormatRow`:F$lastFormatRow").Copy()
    } else {
        $ws.Range("A$normalFormatRow`:F$normalFormatRow").Copy()
    }
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $d[1]

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $d[2]

    $ws.Cells.Item($r, 3).Value = $d[3]

    if ($d[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $d[4] }
    if ($d[5] -ne $null) { $ws.Cells.Item($r, 5).Value = $d[5] }
    if ($d[6] -ne $null) { $ws.Cells.Item($r, 6).Value = $d[6] }
}

$excel.CutCopyMode = 0

# Restore the view state: scrolled so row 109 is at the top, with the cursor
# left on H162 (below/right of the pasted data), matching where editing ended.
$ws.Activate()
$ws.Range("H162").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 109
$win.ScrollColumn = 1
